# Edit script for event_category.xlsx
# 1) Add 19 new event-category rows (430-448) with a new monospace style on col A
# 2) Add 19 trailing empty styled cells in col B (505-523)
# 3) Break the hidden _FilterDatabase defined name ( -> #REF! )
# 4) Move the active selection to D9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    ,@(430, "FLOOD/RAIN/WINDS", "TIDE / FLOOD / SURF")
    ,@(431, "COOL AND WET", "COLD")
    ,@(432, "COLD AND WET CONDITIONS", "COLD")
    ,@(433, "EXCESSIVE WETNESS", "RAIN")
    ,@(434, "SMALL STREAM FLOOD", "TIDE / FLOOD / SURF")
    ,@(435, "HVY RAIN", "RAIN")
    ,@(436, "HAIL 150", "HAIL")
    ,@(437, "HAIL 075", "HAIL")
    ,@(438, "HAIL 125", "HAIL")
    ,@(439, "THUNDERSTORM WIND G60", "THUNDERSTORM/LIGHTNING")
    ,@(440, "THUNDERSTORM WINDS G60", "THUNDERSTORM/LIGHTNING")
    ,@(441, "HARD FREEZE", "COLD")
    ,@(442, "HAIL 200", "HAIL")
    ,@(443, "THUNDERSTORM WIND.", "THUNDERSTORM/LIGHTNING")
    ,@(444, "TORNADOES", "TORNADO")
    ,@(445, "Unseasonable Cold", "COLD")
    ,@(446, "Early Frost", "COLD")
    ,@(447, "AGRICULTURAL FREEZE", "COLD")
    ,@(448, "UNSEASONAL RAIN", "RAIN")
)

# Build the new cell style (Lucida Console 10pt black, vertically centered) once on a
# scratch cell, then fan it out with PasteSpecial so the engine only ever mints the
# font/xf records a single time instead of once per target cell.
$scratch = $ws.Cells.Item(2000, 1)
$scratch.Value = "x"
$scratch.Font.Name = "Lucida Console"
$scratch.Font.Size = 10
$scratch.Font.Color = 0
$scratch.Font.Family = 3
$scratch.VerticalAlignment = -4108
$scratch.Copy() | Out-Null

foreach ($r in $newRows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 1).PasteSpecial(-4122) | Out-Null
}

for ($i = 505; $i -le 523; $i++) {
    $ws.Cells.Item($i, 2).PasteSpecial(-4122) | Out-Null
}

$scratch.Clear() | Out-Null

# The filtered range grew stale/was removed -> Excel leaves a broken reference behind.
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=event_category!#REF!"

# Leave the cursor where the author left it
$ws.Range("D9").Select() | Out-Null
